$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Fagnant2024" sheet after the current last sheet (Bajgain2020)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Fagnant2024"

# ---------------------------------------------------------------------------
# 2. Fill in the sheet contents
#
# NOTE: new shared-string entries are interned in the order the *strings*
# are first written. To reproduce the exact <sst> ordering of the target
# workbook, the string-introducing cells are written first, in this order:
#   A1, A2, A7, K6, K19, B7, K7
# before the remaining (index-reusing / numeric) cells are filled in.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value  = "Taken from fig 5 in pub"
$ws.Cells.Item(2,1).Value  = "Took all agronomic managements, middle of box distribution (median)"
$ws.Cells.Item(7,1).Value  = "Fagnant et al. 2024"
$ws.Cells.Item(6,11).Value = "notes"
$ws.Cells.Item(19,11).Value = "split application"
$ws.Cells.Item(7,2).Value  = "Belgium"
$ws.Cells.Item(7,11).Value = "Planting density was 250 seeds m-2 while actual plant density was 176 seed m-2"

# header row (row 6)
$ws.Cells.Item(6,1).Value  = "citation"
$ws.Cells.Item(6,2).Value  = "site"
$ws.Cells.Item(6,3).Value  = "crop"
$ws.Cells.Item(6,4).Value  = "year"
$ws.Cells.Item(6,5).Value  = "nfert_kgha"
$ws.Cells.Item(6,6).Value  = "grain_kgha"
$ws.Cells.Item(6,7).Value  = "biomass_kgha"
$ws.Cells.Item(6,8).Value  = "hi_pct"
$ws.Cells.Item(6,9).Value  = "pd_seedsm2"
$ws.Cells.Item(6,10).Value = "rowsp_cm"

# row 7
$ws.Cells.Item(7,3).Value  = "IWG"
$ws.Cells.Item(7,4).Value  = 1
$ws.Cells.Item(7,5).Value  = 0
$ws.Cells.Item(7,6).Value  = 590
$ws.Cells.Item(7,9).Value  = 250
$ws.Cells.Item(7,10).Value = 25

# row 8
$ws.Cells.Item(8,4).Value = 2
$ws.Cells.Item(8,6).Value = 1050

# row 9
$ws.Cells.Item(9,4).Value = 3
$ws.Cells.Item(9,6).Value = 496

# row 10
$ws.Cells.Item(10,4).Value = 4
$ws.Cells.Item(10,6).Value = 843

# row 11
$ws.Cells.Item(11,4).Value = 1
$ws.Cells.Item(11,5).Value = 50
$ws.Cells.Item(11,6).Value = 845

# row 12
$ws.Cells.Item(12,4).Value = 2
$ws.Cells.Item(12,6).Value = 903

# row 13
$ws.Cells.Item(13,4).Value = 3
$ws.Cells.Item(13,6).Value = 700

# row 14
$ws.Cells.Item(14,4).Value = 4
$ws.Cells.Item(14,6).Value = 1092

# row 15
$ws.Cells.Item(15,4).Value = 1
$ws.Cells.Item(15,5).Value = 100
$ws.Cells.Item(15,6).Value = 990

# row 16
$ws.Cells.Item(16,4).Value = 2
$ws.Cells.Item(16,6).Value = 896

# row 17
$ws.Cells.Item(17,4).Value = 3
$ws.Cells.Item(17,6).Value = 843

# row 18
$ws.Cells.Item(18,4).Value = 4
$ws.Cells.Item(18,6).Value = 1053

# row 19
$ws.Cells.Item(19,4).Value  = 1
$ws.Cells.Item(19,5).Value  = 100
$ws.Cells.Item(19,6).Value  = 1088

# row 20
$ws.Cells.Item(20,4).Value  = 2
$ws.Cells.Item(20,6).Value  = 798
$ws.Cells.Item(20,11).Value = "split application"

# row 21
$ws.Cells.Item(21,4).Value  = 3
$ws.Cells.Item(21,6).Value  = 871
$ws.Cells.Item(21,11).Value = "split application"

# row 22
$ws.Cells.Item(22,4).Value  = 4
$ws.Cells.Item(22,6).Value  = 1022
$ws.Cells.Item(22,11).Value = "split application"

# ---------------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping
#    (Fernandez2020's selection is set first; Fagnant2024 is activated last
#    so it ends up as the truly-active / tabSelected sheet.)
# ---------------------------------------------------------------------------
$fernandez = $wb.Worksheets.Item("Fernandez2020")
[void]$fernandez.Range("A1:K16").Select()

[void]$ws.Activate()
[void]$ws.Range("G18").Select()
